$wb = $excel.ActiveWorkbook
$wsShort = $wb.Worksheets.Item("Short Term")
$wsMedium = $wb.Worksheets.Item("Medium Term")

# --- Update existing cell values (revised figures) ---

# Short Term sheet
$wsShort.Range("F110").Value = 21.82
$wsShort.Range("F115").Value = -8.73
$wsShort.Range("G115").Value = -6.23
$wsShort.Range("F116").Value = 16.59
$wsShort.Range("G116").Value = 11.96
$wsShort.Range("F117").Value = 7.95
$wsShort.Range("G117").Value = 8.06
$wsShort.Range("F118").Value = 40.69
$wsShort.Range("G118").Value = -4.7
$wsShort.Range("F119").Value = 21.06
$wsShort.Range("G119").Value = 1.66
$wsShort.Range("F120").Value = 28.57
$wsShort.Range("G120").Value = -17.12
$wsShort.Range("B122").Value = 31.07
$wsShort.Range("C122").Value = 31.63
$wsShort.Range("D122").Value = 1.45
$wsShort.Range("E122").Value = 34.52
$wsShort.Range("F122").Value = 31.86
$wsShort.Range("G122").Value = 13.87
$wsShort.Range("B123").Value = 3.73
$wsShort.Range("C123").Value = 9.2
$wsShort.Range("D123").Value = 14.08
$wsShort.Range("E123").Value = 66.13
$wsShort.Range("F123").Value = 56.75
$wsShort.Range("G123").Value = 17.22
$wsShort.Range("B124").Value = 36.54
$wsShort.Range("C124").Value = 37.83
$wsShort.Range("D124").Value = -4.34
$wsShort.Range("E124").Value = 102.19
$wsShort.Range("F124").Value = 95.92
$wsShort.Range("G124").Value = 10.76
$wsShort.Range("B125").Value = -36.81
$wsShort.Range("C125").Value = -38.2
$wsShort.Range("D125").Value = -12.31
$wsShort.Range("B126").Value = 11.31
$wsShort.Range("C126").Value = 6.98
$wsShort.Range("D126").Value = 2.13
$wsShort.Range("B127").Value = -23.75
$wsShort.Range("C127").Value = -28.18
$wsShort.Range("D127").Value = 6.06
$wsShort.Range("E127").Value = 10.68
$wsShort.Range("F127").Value = 3.2
$wsShort.Range("G127").Value = 4.52

# Medium Term sheet
$wsMedium.Range("B96").Value = 1.65
$wsMedium.Range("C96").Value = -5.86
$wsMedium.Range("D96").Value = -4.22
$wsMedium.Range("B97").Value = 2.26
$wsMedium.Range("C97").Value = -3.56
$wsMedium.Range("D97").Value = -4.98
$wsMedium.Range("B98").Value = 5.96
$wsMedium.Range("C98").Value = 0.4
$wsMedium.Range("D98").Value = -4.07
$wsMedium.Range("C99").Value = 4.76
$wsMedium.Range("D99").Value = -1.05
$wsMedium.Range("C100").Value = 7.33
$wsMedium.Range("D100").Value = 0.29
$wsMedium.Range("B101").Value = 8.13
$wsMedium.Range("C101").Value = 7.04
$wsMedium.Range("D101").Value = -0.81
$wsMedium.Range("B102").Value = 4.8
$wsMedium.Range("C102").Value = 6.25
$wsMedium.Range("D102").Value = 0.1
$wsMedium.Range("B103").Value = 4.53
$wsMedium.Range("C103").Value = 8.33
$wsMedium.Range("D103").Value = 2.36
$wsMedium.Range("B104").Value = 21.69
$wsMedium.Range("C104").Value = 14.73
$wsMedium.Range("D104").Value = 7.45
$wsMedium.Range("B105").Value = 23.08
$wsMedium.Range("C105").Value = 13.84
$wsMedium.Range("D105").Value = 9.33
$wsMedium.Range("B106").Value = 29.79
$wsMedium.Range("C106").Value = 17.08
$wsMedium.Range("D106").Value = 12.24
$wsMedium.Range("B107").Value = 21.61
$wsMedium.Range("C107").Value = 21.65
$wsMedium.Range("D107").Value = 14.12
$wsMedium.Range("B108").Value = 25.78
$wsMedium.Range("C108").Value = 24.42
$wsMedium.Range("D108").Value = 15.2
$wsMedium.Range("B109").Value = 35.01
$wsMedium.Range("C109").Value = 32.38
$wsMedium.Range("D109").Value = 20.2
$wsMedium.Range("B110").Value = 61.95
$wsMedium.Range("C110").Value = 43.1
$wsMedium.Range("D110").Value = 29.18
$wsMedium.Range("B111").Value = 54.83
$wsMedium.Range("C111").Value = 41.12
$wsMedium.Range("D111").Value = 27.69
$wsMedium.Range("B112").Value = 46.88
$wsMedium.Range("C112").Value = 41.4
$wsMedium.Range("D112").Value = 29.58
$wsMedium.Range("B113").Value = 16.03
$wsMedium.Range("C113").Value = 38.94
$wsMedium.Range("D113").Value = 30.86

# --- Append new monthly rows (Jul/Aug/Sep 2025), copying date formatting from the row above ---

# Short Term sheet - new rows
$wsShort.Range("A127").Copy()
$wsShort.Range("A128").PasteSpecial(-4122)
$wsShort.Range("A128").Value = 45839
$wsShort.Range("B128").Value = 2.08
$wsShort.Range("C128").Value = 3.66
$wsShort.Range("D128").Value = 5.62
$wsShort.Range("E128").Value = 1.8
$wsShort.Range("F128").Value = -5.87
$wsShort.Range("G128").Value = -0.19

$wsShort.Range("A128").Copy()
$wsShort.Range("A129").PasteSpecial(-4122)
$wsShort.Range("A129").Value = 45870
$wsShort.Range("B129").Value = 1.27
$wsShort.Range("C129").Value = -0.04
$wsShort.Range("D129").Value = 3.21
$wsShort.Range("E129").Value = -0.48
$wsShort.Range("F129").Value = -6.92
$wsShort.Range("G129").Value = 2.83

$wsShort.Range("A129").Copy()
$wsShort.Range("A130").PasteSpecial(-4122)
$wsShort.Range("A130").Value = 45901
$wsShort.Range("B130").Value = 50.63
$wsShort.Range("C130").Value = 57.08
$wsShort.Range("D130").Value = -11.87
$wsShort.Range("E130").Value = 39.34
$wsShort.Range("F130").Value = 27.81
$wsShort.Range("G130").Value = 3.72


# Medium Term sheet - new rows
$wsMedium.Range("A113").Copy()
$wsMedium.Range("A114").PasteSpecial(-4122)
$wsMedium.Range("A114").Value = 45839
$wsMedium.Range("B114").Value = 9.73
$wsMedium.Range("C114").Value = 32.6
$wsMedium.Range("D114").Value = 28.69

$wsMedium.Range("A114").Copy()
$wsMedium.Range("A115").PasteSpecial(-4122)
$wsMedium.Range("A115").Value = 45870
$wsMedium.Range("B115").Value = -3.32
$wsMedium.Range("C115").Value = 22.76
$wsMedium.Range("D115").Value = 27.32

$wsMedium.Range("A115").Copy()
$wsMedium.Range("A116").PasteSpecial(-4122)
$wsMedium.Range("A116").Value = 45901
$wsMedium.Range("B116").Value = 6.78
$wsMedium.Range("C116").Value = 11.25
$wsMedium.Range("D116").Value = 26.37
